$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 10-16: relabel (column B) and set recomputed averaged-intensity values (C:M) ---
# The averaging-scheme catalogue was reordered/extended (new Gaussian-Quadrature + 3 Spiral schemes
# inserted ahead of the existing NoRotation/Rotation/HexGrid schemes), so every row from 10 down shifts
# which scheme it reports and gets freshly computed numbers.

# Row 10: Gaussian-Quadrature
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9997970327128718
$ws.Range("D10").Value = 0.9697660149721196
$ws.Range("E10").Value = 0.9997385364559017
$ws.Range("F10").Value = 0.9997970327128718
$ws.Range("G10").Value = 0.9699415614231334
$ws.Range("H10").Value = 1.001031999558077
$ws.Range("I10").Value = 0.9941468951873086
$ws.Range("J10").Value = 0.9697660149721196
$ws.Range("K10").Value = 0.9847522757140106
$ws.Range("L10").Value = 0.9922746542134413
$ws.Range("M10").Value = 0.9890703400515687

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9921906771602967
$ws.Range("D11").Value = 0.9798195914221737
$ws.Range("E11").Value = 0.9947802115617808
$ws.Range("F11").Value = 0.9921906771602967
$ws.Range("G11").Value = 0.9835719285697239
$ws.Range("H11").Value = 1.000151500738855
$ws.Range("I11").Value = 0.9941399773397617
$ws.Range("J11").Value = 0.9798195914221737
$ws.Range("K11").Value = 0.9872999014919772
$ws.Range("L11").Value = 0.9897452893261369
$ws.Range("M11").Value = 0.9907756477987654

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9921551576612165
$ws.Range("D12").Value = 0.9798645986460404
$ws.Range("E12").Value = 0.9947869749823363
$ws.Range("F12").Value = 0.9921551576612165
$ws.Range("G12").Value = 0.9836799678045381
$ws.Range("H12").Value = 1.000159340623365
$ws.Range("I12").Value = 0.9941451741110581
$ws.Range("J12").Value = 0.9798645986460404
$ws.Range("K12").Value = 0.9873257868141883
$ws.Range("L12").Value = 0.9897404722377023
$ws.Range("M12").Value = 0.9907985356380923

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9921946967321978
$ws.Range("D13").Value = 0.9797980644512084
$ws.Range("E13").Value = 0.9948272211591419
$ws.Range("F13").Value = 0.9921946967321978
$ws.Range("G13").Value = 0.9836395797508477
$ws.Range("H13").Value = 1.000160864735405
$ws.Range("I13").Value = 0.994161525821143
$ws.Range("J13").Value = 0.9797980644512084
$ws.Range("K13").Value = 0.9873126428051752
$ws.Range("L13").Value = 0.9897536697686865
$ws.Range("M13").Value = 0.9907969921083241

# Row 14: NoRotation-tilt60deg
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.9968679999999996
$ws.Range("D14").Value = 0.9314519999999998
$ws.Range("E14").Value = 0.9981239999999998
$ws.Range("F14").Value = 0.9968679999999996
$ws.Range("G14").Value = 0.9341439999999993
$ws.Range("H14").Value = 1.053979999999997
$ws.Range("I14").Value = 0.9972920000000007
$ws.Range("J14").Value = 0.9314519999999998
$ws.Range("K14").Value = 0.9647879999999998
$ws.Range("L14").Value = 0.9808279999999997
$ws.Range("M14").Value = 0.9853099999999994

# Row 15: Rotation-NoTilt
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0.89
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.89
$ws.Range("H15").Value = 1.1
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 0.89
$ws.Range("K15").Value = 0.9450000000000001
$ws.Range("L15").Value = 0.9725
$ws.Range("M15").Value = 0.9800000000000001

# Row 16: Rotation-60detTilt
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9968750080000031
$ws.Range("D16").Value = 0.9313363466239953
$ws.Range("E16").Value = 0.9981250047999984
$ws.Range("F16").Value = 0.9968750080000031
$ws.Range("G16").Value = 0.9343060402176037
$ws.Range("H16").Value = 1.053958451199999
$ws.Range("I16").Value = 0.9972916736000018
$ws.Range("J16").Value = 0.9313363466239953
$ws.Range("K16").Value = 0.9647306757119969
$ws.Range("L16").Value = 0.980802841856
$ws.Range("M16").Value = 0.9853154207402669

# --- Append 3 new rows (17-19) for the new spiral/HexGrid averaging schemes ---

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A16:M16").Copy($ws.Range("A17"))
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9902816289678632
$ws.Range("D17").Value = 0.9916090644196027
$ws.Range("E17").Value = 0.9911008089428689
$ws.Range("F17").Value = 0.9902816289678632
$ws.Range("G17").Value = 0.9905680171969348
$ws.Range("H17").Value = 0.9900085496758524
$ws.Range("I17").Value = 0.9912251169154227
$ws.Range("J17").Value = 0.9916090644196027
$ws.Range("K17").Value = 0.9913549366812358
$ws.Range("L17").Value = 0.9908182828245495
$ws.Range("M17").Value = 0.9907988643530908

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A17:M17").Copy($ws.Range("A18"))
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9905469845476249
$ws.Range("D18").Value = 0.9969328779608267
$ws.Range("E18").Value = 0.9910885936214049
$ws.Range("F18").Value = 0.9905469845476249
$ws.Range("G18").Value = 0.9923287446189193
$ws.Range("H18").Value = 0.9846493877629974
$ws.Range("I18").Value = 0.9902983081668018
$ws.Range("J18").Value = 0.9969328779608267
$ws.Range("K18").Value = 0.9940107357911159
$ws.Range("L18").Value = 0.9922788601693704
$ws.Range("M18").Value = 0.990974149446429

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A18:M18").Copy($ws.Range("A19"))
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9888161938134855
$ws.Range("D19").Value = 1.004325867851453
$ws.Range("E19").Value = 0.9886897954349229
$ws.Range("F19").Value = 0.9888161938134855
$ws.Range("G19").Value = 0.9993359301377945
$ws.Range("H19").Value = 0.9801588706975204
$ws.Range("I19").Value = 0.9886335434724505
$ws.Range("J19").Value = 1.004325867851453
$ws.Range("K19").Value = 0.9965078316431879
$ws.Range("L19").Value = 0.9926620127283368
$ws.Range("M19").Value = 0.9916600335679377

